$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": a new salesperson ("GUERRERO GARCIA OLIMPIA
# ANNABELLE") is inserted into the OFICINA-CATAECSA group, alphabetically
# right before "JAIME COELLO ALBERTO FERNANDO" (row 304). Inserting a whole
# row there pushes every following row (304-358) down by one, which is
# exactly what the diff shows (each name/value moves from row N to N+1).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(304).Insert()

$ws1.Cells.Item(304, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(304, 2).Value = "GUERRERO GARCIA OLIMPIA ANNABELLE"
for ($c = 3; $c -le 18; $c++) {
  $ws1.Cells.Item(304, $c).Value = 0
}

# The summary/footer row (now shifted from row 358 to row 359) shows counts
# like "4 de 356" -- the total number of salespeople grew from 356 to 357,
# so every label in that row needs the count bumped.
for ($c = 3; $c -le 18; $c++) {
  $cell = $ws1.Cells.Item(359, $c)
  $text = $cell.Text
  $cell.Value = $text.Replace("356", "357")
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same new salesperson inserted in the same spot
# (before "JAIME COELLO ALBERTO FERNANDO", row 308 here), shifting rows
# 308-362 down by one.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(308).Insert()

$ws2.Cells.Item(308, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(308, 2).Value = "GUERRERO GARCIA OLIMPIA ANNABELLE"
for ($c = 3; $c -le 7; $c++) {
  $ws2.Cells.Item(308, $c).Value = 0
}
